# LocaDemo.xlsx update — "Improved LocaExcelBridge All sheets supported"
#
# Adds a second localisation sheet ("Tabelle2") holding a "cat" row built the
# same way as the existing "dog" row on "Tabelle1", and updates the on-sheet
# selection state left behind by the edit.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Add the new worksheet after the existing one and name it "Tabelle2"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Tabelle2"

# ---------------------------------------------------------------------
# 2) Column widths for the description columns (D, E) — ~36.29 / ~30.57
#    characters, matching the widened "Katzenbeschreibung" columns
# ---------------------------------------------------------------------
$ws2.Columns.Item(4).ColumnWidth = 36.28515625
$ws2.Columns.Item(5).ColumnWidth = 30.5703125


# ---------------------------------------------------------------------
# 3) Header rows — identical layout/wording to Tabelle1's header rows
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Laber"
$ws2.Range("B1").Value = "Bla"
$ws2.Range("C1").Value = "Dont use"
$ws2.Range("D1").Value = "this"
$ws2.Range("E1").Value = "row"
$ws2.Range("F1").Value = "Sing."
$ws2.Range("G1").Value = "Plur."
$ws2.Range("H1").Value = "Sing."
$ws2.Range("I1").Value = "Plur."

$ws2.Range("A2").Value = "key"
$ws2.Range("B2").Value = "en"
$ws2.Range("C2").Value = "de"
$ws2.Range("D2").Value = "de"
$ws2.Range("E2").Value = "de"
$ws2.Range("F2").Value = "en"
$ws2.Range("G2").Value = "en"
$ws2.Range("H2").Value = "de"
$ws2.Range("I2").Value = "de"

# ---------------------------------------------------------------------
# 4) Data row — the new "cat" localisation entry
# ---------------------------------------------------------------------
$ws2.Range("A3").Value = "excel_cat"
$ws2.Range("B3").Value = "cat"
$ws2.Range("C3").Value = "Katze"
$ws2.Range("D3").Value = "Katzenbeschreibung"
$ws2.Range("E3").Value = "Noch ne Katzenbeschreibung"
$ws2.Range("F3").Value = "cat"
$ws2.Range("G3").Value = "cats"
$ws2.Range("H3").Value = "Katze"
$ws2.Range("I3").Value = "Katzen"

# ---------------------------------------------------------------------
# 5) Page margins to match the workbook's existing (2cm top/bottom) style
# ---------------------------------------------------------------------
$ws2.PageSetup.TopMargin = 56.692913399999995
$ws2.PageSetup.BottomMargin = 56.692913399999995

# ---------------------------------------------------------------------
# 6) Selection / view state left behind on each sheet
# ---------------------------------------------------------------------
$ws1.Range("A1:I4").Select()

$ws2.Range("B3").Select()
$ws2.Activate()
